$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 2 (pushes existing rows 2-13 down to 3-14),
# so the table gains the new "코칩" record at the top of the dataset.
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "'2024-04-15"
$ws.Range("B2").Value = "'2024-04-19"
$ws.Range("C2").Value = "'2024-05-07"
$ws.Range("D2").Value = "한국"
$ws.Range("E2").Value = "코칩"
$ws.Range("F2").Value = 1500000
$ws.Range("G2").Value = 1500000
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 11000
$ws.Range("J2").Value = 14000
$ws.Range("K2").Value = 8503460
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 18000
$ws.Range("N2").Value = "988.32:1"
$ws.Range("O2").Value = "'13.19%"
$ws.Range("P2").Value = 47284698907
$ws.Range("Q2").Value = 38750429966
$ws.Range("R2").Value = 25900014771
$ws.Range("S2").Value = 7595091433
$ws.Range("T2").Value = 5807002440
$ws.Range("U2").Value = 3668321605
$ws.Range("V2").Value = 5701880294
$ws.Range("W2").Value = 4780312126
$ws.Range("X2").Value = 4195570793
$ws.Range("Y2").Value = "소형 및 초소형 슈퍼커패시터"

# The new row should have no special style/border (same look as the other
# plain data rows), so strip whatever formatting Insert()/quote-prefix left.
$ws.Rows.Item(2).ClearFormats()

# The table only ever shows 12 data rows, so the oldest record (which is
# now pushed to row 14: 케이엔알시스템 / DB, NH) drops off the bottom.
$ws.Rows.Item(14).Delete()
